$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date" for dbf4aa4f...md (row 5)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-10-20 08:31:21"

# zh-cn sheet: column H = "Latest Handoff Datetime" for dbf4aa4f...md (row 5)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-10-20 08:31:08"

# de-de sheet: column H = "Latest Handoff Datetime" for dbf4aa4f...md (row 5)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-10-20 08:31:21"
